$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.433.51"
$ws.Range("E2").Value = "  +1.99%  "

$ws.Range("D3").Value = "3.579.42"
$ws.Range("E3").Value = "  +0.42%  "

$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").Value = "207.04"
$ws.Range("E5").Value = "  +8.56%  "

$ws.Range("D6").Value = "561.60"
$ws.Range("E6").Value = "  -1.30%  "

$ws.Range("D7").Value = "0.608"
$ws.Range("E7").Value = "  -0.99%  "

$ws.Range("E8").Value = "  -0.18%  "

$ws.Range("D9").Value = "0.673"
$ws.Range("E9").Value = "  -0.19%  "

$ws.Range("D10").Value = "63.34"
$ws.Range("E10").Value = "  +13.29%  "

$ws.Range("D11").Value = "0.146"
$ws.Range("E11").Value = "  -1.88%  "

$ws.Range("D12").Value = "0.0000279"
$ws.Range("E12").Value = "  +3.32%  "

$ws.Range("D13").Value = "10.10"
$ws.Range("E13").Value = "  +2.74%  "

$ws.Range("D14").Value = "4.131.71"
$ws.Range("E14").Value = "  -0.02%  "

$ws.Range("D15").Value = "3.566.82"
$ws.Range("E15").Value = "  -0.05%  "

$ws.Range("E16").Value = "  +0.05%  "

$ws.Range("D17").Value = "19.04"
$ws.Range("E17").Value = "  +4.88%  "

$ws.Range("D18").Value = "68.058.58"
$ws.Range("E18").Value = "  +1.72%  "

$ws.Range("D19").Value = "12.15"
$ws.Range("E19").Value = "  -0.08%  "

$ws.Range("D20").Value = "1.06"
$ws.Range("E20").Value = "  +0.14%  "

$ws.Range("D21").Value = "397.69"
$ws.Range("E21").Value = "  -0.63%  "

$ws.Range("B22").Value = "RenderToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D22").Value = "12.32"
$ws.Range("E22").Value = "  +4.23%  "

$ws.Range("B23").Value = "PancakeSwap"
$ws.Range("C23").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D23").Value = "4.11"
$ws.Range("E23").Value = "  -0.65%  "

$ws.Range("D24").Value = "84.13"
$ws.Range("E24").Value = "  -1.77%  "

$ws.Range("B25").Value = "ImmutableX"
$ws.Range("C25").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D25").Value = "2.87"
$ws.Range("E25").Value = "  -1.38%  "

$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").Value = "12.39"
$ws.Range("E26").Value = "  -0.59%  "

$ws.Range("D27").Value = "3.86"
$ws.Range("E27").Value = "  +5.58%  "

$ws.Range("D28").Value = "9.09"
$ws.Range("E28").Value = "  +1.38%  "

$ws.Range("D29").Value = "717.20"
$ws.Range("E29").Value = "  +11.81%  "

$ws.Range("D30").Value = "31.35"
$ws.Range("E30").Value = "  +0.68%  "

$ws.Range("D31").Value = "7.52"
$ws.Range("E31").Value = "  -3.56%  "

$ws.Range("D32").Value = "12.02"
$ws.Range("E32").Value = "  -0.67%  "

$ws.Range("D33").Value = "63.73"
$ws.Range("E33").Value = "  -0.02%  "

$ws.Range("D34").Value = "0.112"
$ws.Range("E34").Value = "  -1.70%  "

$ws.Range("D35").Value = "41.19"
$ws.Range("E35").Value = "  -2.56%  "

$ws.Range("D36").Value = "0.419"
$ws.Range("E36").Value = "  +3.56%  "

$ws.Range("E37").Value = "  -0.17%  "

$ws.Range("D38").Value = "3.22"
$ws.Range("E38").Value = "  +6.70%  "

$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").Value = "3.13"
$ws.Range("E39").Value = "  +28.77%  "

$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "3.140.60"
$ws.Range("E40").Value = "  -1.22%  "

$ws.Range("D41").Value = "0.0₃0730"
$ws.Range("E41").Value = "  -4.33%  "

$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "0.131"
$ws.Range("E42").Value = "  -1.44%  "

$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "0.995"
$ws.Range("E43").Value = "  -0.35%  "

$ws.Range("D44").Value = "2.58"
$ws.Range("E44").Value = "  -4.37%  "

$ws.Range("E45").Value = "  +9.20%  "

$ws.Range("D46").Value = "0.0411"
$ws.Range("E46").Value = "  -0.50%  "

$ws.Range("D47").Value = "0.130"
$ws.Range("E47").Value = "  -0.08%  "

$ws.Range("D48").Value = "3.05"
$ws.Range("E48").Value = "  -1.00%  "

$ws.Range("D49").Value = "8.61"
$ws.Range("E49").Value = "  +1.31%  "

$ws.Range("D50").Value = "138.46"
$ws.Range("E50").Value = "  -2.15%  "

$ws.Range("D51").Value = "2.69"
$ws.Range("E51").Value = "  -0.18%  "
